# Insert a new weekly data row above current row 7 (pushing existing rows 7-28 down to 8-29),
# then populate the new row with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; existing row 7 and below shift down to 8.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new record's values.
$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44481
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 100112022
$ws.Cells.Item(7, 7).Value = "Arveja Verde"
$ws.Cells.Item(7, 8).Value = "Perfection"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 63
$ws.Cells.Item(7, 11).Value = 22000
$ws.Cells.Item(7, 12).Value = 23000
$ws.Cells.Item(7, 13).Value = 22476
$ws.Cells.Item(7, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 899
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"
